$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 8).Value = [double]"0.1055113828966178"

$ws.Cells.Item(3, 2).Value = [double]"0.003077050901595573"
$ws.Cells.Item(3, 3).Value = [double]"0.0007090834989552793"
$ws.Cells.Item(3, 4).Value = [double]"5.155846636029093"
$ws.Cells.Item(3, 5).Value = [double]"0.08345614489873396"
$ws.Cells.Item(3, 6).Value = [double]"0.00168726747658944"
$ws.Cells.Item(3, 7).Value = [double]"0.004466834326601707"
$ws.Cells.Item(3, 8).Value = [double]"0.1085884337982134"

$ws.Cells.Item(4, 2).Value = [double]"0.006183902575172915"
$ws.Cells.Item(4, 3).Value = [double]"0.001623603099642614"
$ws.Cells.Item(4, 4).Value = [double]"6.058520379474212"
$ws.Cells.Item(4, 5).Value = [double]"0.03470057094788833"
$ws.Cells.Item(4, 6).Value = [double]"0.003001688329325818"
$ws.Cells.Item(4, 7).Value = [double]"0.009366116821020011"
$ws.Cells.Item(4, 8).Value = [double]"0.1116952854717907"

$ws.Cells.Item(5, 2).Value = [double]"0.01037695341129617"
$ws.Cells.Item(5, 3).Value = [double]"0.002373578590677842"
$ws.Cells.Item(5, 4).Value = [double]"6.858245491748216"
$ws.Cells.Item(5, 5).Value = [double]"0.0002348453347379148"
$ws.Cells.Item(5, 6).Value = [double]"0.005724810501236933"
$ws.Cells.Item(5, 7).Value = [double]"0.0150290963213554"
$ws.Cells.Item(5, 8).Value = [double]"0.115888336307914"

$ws.Cells.Item(6, 2).Value = [double]"0.007057157620615919"
$ws.Cells.Item(6, 3).Value = [double]"0.004298490341215918"
$ws.Cells.Item(6, 4).Value = [double]"2.344374205507154"
$ws.Cells.Item(6, 5).Value = [double]"0.07419427420557795"
$ws.Cells.Item(6, 6).Value = [double]"-0.001367753526678727"
$ws.Cells.Item(6, 7).Value = [double]"0.01548206876791056"
$ws.Cells.Item(6, 8).Value = [double]"0.1125685405172337"

$ws.Cells.Item(7, 2).Value = [double]"0.0085021688254235"
$ws.Cells.Item(7, 3).Value = [double]"0.004310308172033145"
$ws.Cells.Item(7, 4).Value = [double]"1.233395711586459"
$ws.Cells.Item(7, 5).Value = [double]"0.1746252399387233"
$ws.Cells.Item(7, 6).Value = [double]"5.409493788441077e-05"
$ws.Cells.Item(7, 7).Value = [double]"0.01695024271296259"
$ws.Cells.Item(7, 8).Value = [double]"0.1140135517220413"

$ws.Cells.Item(8, 2).Value = [double]"0.009890863204363233"
$ws.Cells.Item(8, 3).Value = [double]"0.005394945667408544"
$ws.Cells.Item(8, 4).Value = [double]"0.8894896090888758"
$ws.Cells.Item(8, 5).Value = [double]"0.1525084530008978"
$ws.Cells.Item(8, 6).Value = [double]"-0.000683065177812207"
$ws.Cells.Item(8, 7).Value = [double]"0.02046479158653867"
$ws.Cells.Item(8, 8).Value = [double]"0.1154022461009811"

$ws.Cells.Item(9, 2).Value = [double]"0.008360746701278858"
$ws.Cells.Item(9, 3).Value = [double]"0.005621032217840073"
$ws.Cells.Item(9, 4).Value = [double]"0.8606591534204926"
$ws.Cells.Item(9, 5).Value = [double]"0.1307706093420732"
$ws.Cells.Item(9, 6).Value = [double]"-0.002656306746397124"
$ws.Cells.Item(9, 7).Value = [double]"0.01937780014895484"
$ws.Cells.Item(9, 8).Value = [double]"0.1138721295978967"

$ws.Cells.Item(10, 2).Value = [double]"-0.1055113828966178"
$ws.Cells.Item(10, 3).Value = [double]"0.0005285358317725725"
$ws.Cells.Item(10, 4).Value = [double]"-229.191924130787"
$ws.Cells.Item(10, 5).Value = [double]"0"
$ws.Cells.Item(10, 6).Value = [double]"-0.1065472978095343"
$ws.Cells.Item(10, 7).Value = [double]"-0.1044754679837013"

$ws.Cells.Item(11, 2).Value = [double]"-0.04978652197816739"
$ws.Cells.Item(11, 3).Value = [double]"0.0005767820399017431"
$ws.Cells.Item(11, 4).Value = [double]"-96.29705670275807"
$ws.Cells.Item(11, 5).Value = [double]"8.81580535912254e-171"
$ws.Cells.Item(11, 6).Value = [double]"-0.05091699807586809"
$ws.Cells.Item(11, 7).Value = [double]"-0.04865604588046671"
$ws.Cells.Item(11, 8).Value = [double]"0.05572486091845043"

$ws.Cells.Item(12, 2).Value = [double]"-0.03913832365960477"
$ws.Cells.Item(12, 3).Value = [double]"0.0005510321437587137"
$ws.Cells.Item(12, 4).Value = [double]"-79.93599789829119"
$ws.Cells.Item(12, 5).Value = [double]"3.409394990344216e-87"
$ws.Cells.Item(12, 6).Value = [double]"-0.04021833070820878"
$ws.Cells.Item(12, 7).Value = [double]"-0.03805831661100076"
$ws.Cells.Item(12, 8).Value = [double]"0.06637305923701306"

$ws.Cells.Item(13, 2).Value = [double]"-0.03617143999336091"
$ws.Cells.Item(13, 3).Value = [double]"0.0005519944846386196"
$ws.Cells.Item(13, 4).Value = [double]"-73.41614956222762"
$ws.Cells.Item(13, 5).Value = [double]"1.429841425579727e-103"
$ws.Cells.Item(13, 6).Value = [double]"-0.03725333318157196"
$ws.Cells.Item(13, 7).Value = [double]"-0.03508954680514985"
$ws.Cells.Item(13, 8).Value = [double]"0.06933994290325692"

$ws.Cells.Item(14, 2).Value = [double]"-0.02985661707435921"
$ws.Cells.Item(14, 3).Value = [double]"0.0005331536511850108"
$ws.Cells.Item(14, 4).Value = [double]"-64.41443175916565"
$ws.Cells.Item(14, 5).Value = [double]"2.665552675886679e-23"
$ws.Cells.Item(14, 6).Value = [double]"-0.03090158278942326"
$ws.Cells.Item(14, 7).Value = [double]"-0.02881165135929518"
$ws.Cells.Item(14, 8).Value = [double]"0.07565476582225861"

$ws.Cells.Item(15, 2).Value = [double]"-0.02692384769151067"
$ws.Cells.Item(15, 3).Value = [double]"0.0005247263668221543"
$ws.Cells.Item(15, 4).Value = [double]"-58.7902041874756"
$ws.Cells.Item(15, 5).Value = [double]"1.176921738007564e-11"
$ws.Cells.Item(15, 6).Value = [double]"-0.02795229617493981"
$ws.Cells.Item(15, 7).Value = [double]"-0.02589539920808152"
$ws.Cells.Item(15, 8).Value = [double]"0.07858753520510715"

$ws.Cells.Item(16, 2).Value = [double]"-0.02531101567749676"
$ws.Cells.Item(16, 3).Value = [double]"0.0005206222631651655"
$ws.Cells.Item(16, 4).Value = [double]"-54.35352884026561"
$ws.Cells.Item(16, 5).Value = [double]"3.181565118844527e-51"
$ws.Cells.Item(16, 6).Value = [double]"-0.02633142025692624"
$ws.Cells.Item(16, 7).Value = [double]"-0.02429061109806729"
$ws.Cells.Item(16, 8).Value = [double]"0.08020036721912106"

$ws.Cells.Item(17, 2).Value = [double]"-0.02423816440304967"
$ws.Cells.Item(17, 3).Value = [double]"0.0005281269613392352"
$ws.Cells.Item(17, 4).Value = [double]"-52.09330666426062"
$ws.Cells.Item(17, 5).Value = [double]"8.358045163770415e-15"
$ws.Cells.Item(17, 6).Value = [double]"-0.02527327794817064"
$ws.Cells.Item(17, 7).Value = [double]"-0.02320305085792872"
$ws.Cells.Item(17, 8).Value = [double]"0.08127321849356815"

$ws.Cells.Item(18, 2).Value = [double]"-0.02067331322706777"
$ws.Cells.Item(18, 3).Value = [double]"0.000530140346063316"
$ws.Cells.Item(18, 4).Value = [double]"-43.63124584275277"
$ws.Cells.Item(18, 5).Value = [double]"6.059164845628032e-06"
$ws.Cells.Item(18, 6).Value = [double]"-0.0217123729597752"
$ws.Cells.Item(18, 7).Value = [double]"-0.01963425349436033"
$ws.Cells.Item(18, 8).Value = [double]"0.08483806966955006"

$ws.Cells.Item(19, 2).Value = [double]"-0.01831499382926514"
$ws.Cells.Item(19, 3).Value = [double]"0.0005282911064811415"
$ws.Cells.Item(19, 4).Value = [double]"-38.00756023052819"
$ws.Cells.Item(19, 5).Value = [double]"6.713332125429022e-06"
$ws.Cells.Item(19, 6).Value = [double]"-0.01935042911587681"
$ws.Cells.Item(19, 7).Value = [double]"-0.01727955854265346"
$ws.Cells.Item(19, 8).Value = [double]"0.08719638906735269"

$ws.Cells.Item(20, 2).Value = [double]"-0.01578134585328985"
$ws.Cells.Item(20, 3).Value = [double]"0.000537749933470827"
$ws.Cells.Item(20, 4).Value = [double]"-30.50240151551372"
$ws.Cells.Item(20, 5).Value = [double]"0.0003942099570098182"
$ws.Cells.Item(20, 6).Value = [double]"-0.01683532016207684"
$ws.Cells.Item(20, 7).Value = [double]"-0.01472737154450287"
$ws.Cells.Item(20, 8).Value = [double]"0.08973003704332796"

$ws.Cells.Item(21, 2).Value = [double]"-0.01052484881736467"
$ws.Cells.Item(21, 3).Value = [double]"0.0005457492115543434"
$ws.Cells.Item(21, 4).Value = [double]"-18.39508621060061"
$ws.Cells.Item(21, 5).Value = [double]"0.03813535985923081"
$ws.Cells.Item(21, 6).Value = [double]"-0.0115945014875032"
$ws.Cells.Item(21, 7).Value = [double]"-0.009455196147226147"
$ws.Cells.Item(21, 8).Value = [double]"0.09498653407925314"

$ws.Cells.Item(22, 2).Value = [double]"-0.008306491239587181"
$ws.Cells.Item(22, 3).Value = [double]"0.0005419461950457809"
$ws.Cells.Item(22, 4).Value = [double]"-12.59158202126708"
$ws.Cells.Item(22, 5).Value = [double]"0.0152044336885292"
$ws.Cells.Item(22, 6).Value = [double]"-0.009368690119423524"
$ws.Cells.Item(22, 7).Value = [double]"-0.007244292359750835"
$ws.Cells.Item(22, 8).Value = [double]"0.09720489165703064"

$ws.Cells.Item(23, 2).Value = [double]"-0.006837935983360133"
$ws.Cells.Item(23, 3).Value = [double]"0.0005516755797213348"
$ws.Cells.Item(23, 4).Value = [double]"-9.971652517880655"
$ws.Cells.Item(23, 5).Value = [double]"0.08000060044236663"
$ws.Cells.Item(23, 6).Value = [double]"-0.007919204214158898"
$ws.Cells.Item(23, 7).Value = [double]"-0.005756667752561367"
$ws.Cells.Item(23, 8).Value = [double]"0.09867344691325769"

$ws.Cells.Item(24, 2).Value = [double]"-0.005598038764819321"
$ws.Cells.Item(24, 3).Value = [double]"0.0005439974083955801"
$ws.Cells.Item(24, 4).Value = [double]"-7.54673023002526"
$ws.Cells.Item(24, 5).Value = [double]"0.0701315123840837"
$ws.Cells.Item(24, 6).Value = [double]"-0.006664257970774811"
$ws.Cells.Item(24, 7).Value = [double]"-0.004531819558863833"
$ws.Cells.Item(24, 8).Value = [double]"0.0999133441317985"

$ws.Cells.Item(25, 2).Value = [double]"-0.004194242488561243"
$ws.Cells.Item(25, 3).Value = [double]"0.0005320503509792121"
$ws.Cells.Item(25, 4).Value = [double]"-5.247201438819127"
$ws.Cells.Item(25, 5).Value = [double]"0.09099229916336422"
$ws.Cells.Item(25, 6).Value = [double]"-0.005237045816010939"
$ws.Cells.Item(25, 7).Value = [double]"-0.003151439161111548"
$ws.Cells.Item(25, 8).Value = [double]"0.1013171404080566"

$ws.Cells.Item(26, 2).Value = [double]"0.01288929695319745"
$ws.Cells.Item(26, 3).Value = [double]"0.001111829825022671"
$ws.Cells.Item(26, 4).Value = [double]"2.244962652974674"
$ws.Cells.Item(26, 5).Value = [double]"0.01017539056856102"
$ws.Cells.Item(26, 6).Value = [double]"0.01071014369447223"
$ws.Cells.Item(26, 7).Value = [double]"0.01506845021192268"
$ws.Cells.Item(26, 8).Value = [double]"0.1184006798498153"
